$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text for B1 (shared string "value" -> "first_release_value")
$ws.Range("B1").Value = "first_release_value"

# New A:B data block (rows 2..84), replacing/extending the old A2:B53 block
$data = @(
    @(38398, 2.129035685702092),
    @(38487, 2.893366588103589),
    @(38579, 2.955627504208366),
    @(38671, 1.167401765354654),
    @(38763, 4.89920306769092),
    @(38852, 1.757302400870174),
    @(38944, 3.311126880990756),
    @(39036, 5.317531793333302),
    @(39128, -0.3124417888986386),
    @(39217, 0.7960900427496256),
    @(39309, 3.053483582409172),
    @(39401, 1.262446517965742),
    @(39493, 2.427463180990628),
    @(39583, -0.2156301208428602),
    @(39675, -0.3737946170534627),
    @(39767, -8.061199066098098),
    @(39859, -9.679894544238294),
    @(39948, -0.8195522748085864),
    @(40040, 3.376547379332465),
    @(40132, 2.274906658561562),
    @(40224, 2.608522136020966),
    @(40313, 7.912074126222677),
    @(40405, 2.327131657487215),
    @(40497, 1.840273951972307),
    @(40589, 2.315788057256768),
    @(40678, 2.755574095503576),
    @(40770, 2.530590017523977),
    @(40862, -1.487154252625388),
    @(40954, 1.726328558848905),
    @(41044, 3.275910395958846),
    @(41136, 1.446992732417414),
    @(41228, -2.418840888438723),
    @(41320, -1.803419669653366),
    @(41409, 2.429152676597084),
    @(41501, 0.07763982310385131),
    @(41593, 2.46356965656642),
    @(41685, 0.1648657400322122),
    @(41774, 1.153911995371715),
    @(41866, 1.946997568649863),
    @(41958, 0.9978838984052771),
    @(42050, 0.8038113817623582),
    @(42139, 1.7958205029814),
    @(42231, 0.2175014152489467),
    @(42323, -0.6421078521012049),
    @(42415, 0.9774607232901644),
    @(42505, 1.244614153452133),
    @(42597, -0.4018974900578201),
    @(42689, 1.682996656587392),
    @(42781, 1.289917584994839),
    @(42870, 1.044515111522799),
    @(42962, 1.660004389679614),
    @(43054, 2.736978867984959),
    @(43146, -1.040686374073175),
    @(43235, 0.6815662815033505),
    @(43327, -0.9),
    @(43419, 0.6525912658484003),
    @(43511, 1.0),
    @(43600, -1.346768218203209),
    @(43692, 1.041951863927608),
    @(43784, -0.2335361222581156),
    @(43876, -3.1),
    @(43966, -20.33630385521235),
    @(44058, 12.42550598425463),
    @(44150, 2.813188392915293),
    @(44242, 1.806909003397877),
    @(44331, 1.087285196410733),
    @(44423, 1.450779288666709),
    @(44515, 0.368631518524424),
    @(44607, -2.021680416268424),
    @(44696, 1.118165634023697),
    @(44788, 0.1289585187160185),
    @(44880, -0.4259770459179748),
    @(44972, -0.6379881427730965),
    @(45061, -0.04993664583679447),
    @(45153, 0.3662234592800075),
    @(45245, -1.023919595764212),
    @(45337, -0.2547854660834332),
    @(45427, 1.565950786385088),
    @(45519, -0.2475555093771362),
    @(45611, 0.386705041387188),
    @(45703, -0.6950112442994083),
    @(45792, -0.3135696110003181),
    @(45884, -0.02110736963342674)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

# Ensure the date-format style used by column A (numFmt, bold font, thin border, centered/top align)
# is applied to the newly added rows in column A (rows 54..84 did not exist before)
$dateRng = $ws.Range("A2:A84")
$dateRng.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$dateRng.Font.Bold = $true
$dateRng.HorizontalAlignment = -4108
$dateRng.VerticalAlignment = -4160
$dateRng.Borders.LineStyle = 1
